$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 525.6667
Write-Host "Done"
